$d = $word.ActiveDocument

# 1. Remove the stray "_GoBack" bookmark (both its start and end markers
#    are removed together when the bookmark itself is deleted).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. After the "氏名" label, add a new run containing "（自署）" so the
#    label reads "氏名（自署）".
$rng = $d.Content
$found = $rng.Find.Execute("氏名", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Collapse(0)
    $rng.InsertAfter("（自署）")
}
